# Update gh-pages output data: increment the "想去人数" (column F) counts
# for several events on the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览": row -> new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 666
$ws1.Range("F9").Value  = 5
$ws1.Range("F10").Value = 10914
$ws1.Range("F13").Value = 294
$ws1.Range("F15").Value = 10712
$ws1.Range("F20").Value = 5359

# Sheet "全部类型": same events, different row numbers -> new F value
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 666
$ws4.Range("F11").Value = 5
$ws4.Range("F13").Value = 10914
$ws4.Range("F16").Value = 294
$ws4.Range("F18").Value = 10712
$ws4.Range("F23").Value = 5359
